$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2336769759450172
$ws.Range("C2").Value = 0.4810996563573883
$ws.Range("J2").Value = 0.01374570446735395
$ws.Range("P2").Value = 0.1855670103092784
$ws.Range("S2").Value = 0.0859106529209622
$ws.Range("B3").Value = 0.006024096385542169
$ws.Range("C3").Value = 0.03614457831325301
$ws.Range("J3").Value = 0.07228915662650602
$ws.Range("P3").Value = 0.6445783132530121
$ws.Range("S3").Value = 0.2409638554216867
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.6216216216216216
$ws.Range("S4").Value = 0.2972972972972973
$ws.Range("B6").Value = 0.03428571428571429
$ws.Range("D6").Value = 0.01714285714285714
$ws.Range("F6").Value = 0.04571428571428571
$ws.Range("J6").Value = 0.2685714285714286
$ws.Range("O6").Value = 0.01142857142857143
$ws.Range("Q6").Value = 0.1257142857142857
$ws.Range("R6").Value = 0.09142857142857143
$ws.Range("S6").Value = 0.4057142857142857
$ws.Range("B7").Value = 0.1256830601092896
$ws.Range("D7").Value = 0.02185792349726776
$ws.Range("F7").Value = 0.0273224043715847
$ws.Range("J7").Value = 0.1256830601092896
$ws.Range("O7").Value = 0.01639344262295082
$ws.Range("Q7").Value = 0.1693989071038251
$ws.Range("R7").Value = 0.08743169398907104
$ws.Range("S7").Value = 0.4262295081967213
$ws.Range("B8").Value = 0.1070707070707071
$ws.Range("D8").Value = 0.0101010101010101
$ws.Range("E8").Value = 0.00202020202020202
$ws.Range("F8").Value = 0.03232323232323232
$ws.Range("J8").Value = 0.1454545454545454
$ws.Range("O8").Value = 0.01414141414141414
$ws.Range("Q8").Value = 0.1232323232323232
$ws.Range("R8").Value = 0.1252525252525253
$ws.Range("S8").Value = 0.4404040404040404
$ws.Range("B9").Value = 0.06626506024096386
$ws.Range("D9").Value = 0.006024096385542169
$ws.Range("E9").Value = 0.006024096385542169
$ws.Range("F9").Value = 0.01204819277108434
$ws.Range("J9").Value = 0.1385542168674699
$ws.Range("O9").Value = 0.01204819277108434
$ws.Range("Q9").Value = 0.1506024096385542
$ws.Range("R9").Value = 0.1746987951807229
$ws.Range("S9").Value = 0.4337349397590362
$ws.Range("B10").Value = 0.09843400447427293
$ws.Range("D10").Value = 0.01789709172259508
$ws.Range("E10").Value = 0.001491424310216256
$ws.Range("F10").Value = 0.05592841163310962
$ws.Range("J10").Value = 0.1588366890380313
$ws.Range("O10").Value = 0.01342281879194631
$ws.Range("Q10").Value = 0.2155108128262491
$ws.Range("R10").Value = 0.09694258016405667
$ws.Range("S10").Value = 0.3415361670395227
$ws.Range("G11").Value = 0.1766561514195584
$ws.Range("J11").Value = 0.08517350157728706
$ws.Range("K11").Value = 0.2618296529968454
$ws.Range("L11").Value = 0.4574132492113565
$ws.Range("S11").Value = 0.01892744479495268
$ws.Range("G12").Value = 0.7073170731707317
$ws.Range("J12").Value = 0.2195121951219512
$ws.Range("L12").Value = 0.03658536585365853
$ws.Range("S12").Value = 0.03658536585365853
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.225
$ws.Range("S13").Value = 0.15
$ws.Range("F15").Value = 0.01595744680851064
$ws.Range("H15").Value = 0.2074468085106383
$ws.Range("I15").Value = 0.03723404255319149
$ws.Range("J15").Value = 0.3563829787234042
$ws.Range("K15").Value = 0.07446808510638298
$ws.Range("M15").Value = 0.01595744680851064
$ws.Range("O15").Value = 0.01595744680851064
$ws.Range("S15").Value = 0.2765957446808511
$ws.Range("F16").Value = 0.01734104046242774
$ws.Range("H16").Value = 0.1907514450867052
$ws.Range("I16").Value = 0.03468208092485549
$ws.Range("J16").Value = 0.3757225433526011
$ws.Range("K16").Value = 0.1849710982658959
$ws.Range("M16").Value = 0.01734104046242774
$ws.Range("N16").Value = 0.005780346820809248
$ws.Range("O16").Value = 0.05202312138728324
$ws.Range("S16").Value = 0.1213872832369942
$ws.Range("F17").Value = 0.0228310502283105
$ws.Range("H17").Value = 0.1986301369863014
$ws.Range("I17").Value = 0.06164383561643835
$ws.Range("J17").Value = 0.4360730593607306
$ws.Range("K17").Value = 0.0867579908675799
$ws.Range("M17").Value = 0.0136986301369863
$ws.Range("O17").Value = 0.0684931506849315
$ws.Range("S17").Value = 0.1118721461187215
$ws.Range("F18").Value = 0.03174603174603174
$ws.Range("H18").Value = 0.1547619047619048
$ws.Range("I18").Value = 0.07539682539682539
$ws.Range("J18").Value = 0.4325396825396826
$ws.Range("K18").Value = 0.1071428571428571
$ws.Range("M18").Value = 0.007936507936507936
$ws.Range("O18").Value = 0.06746031746031746
$ws.Range("S18").Value = 0.123015873015873
$ws.Range("F19").Value = 0.01710730948678071
$ws.Range("H19").Value = 0.2208398133748056
$ws.Range("I19").Value = 0.08320373250388803
$ws.Range("J19").Value = 0.3654743390357698
$ws.Range("K19").Value = 0.1057542768273717
$ws.Range("M19").Value = 0.01866251944012442
$ws.Range("O19").Value = 0.0536547433903577
$ws.Range("S19").Value = 0.135303265940902
